$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "`'247.55"
$ws.Range("D2").Style = "Normal"
$ws.Range("D4").Value = "`'5.460"
$ws.Range("D4").Style = "Normal"
$ws.Range("D6").Value = "`'3.377"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Value = "`'1.031"
$ws.Range("D8").Style = "Normal"
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").Value = "`'0.1439"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "8WazirXWRX"
$ws.Range("B10").Value = "MandalaExchangeToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D10").Value = "`'0.07225"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "9MandalaExchangeTokenMDX"
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").Value = "`'0.03163"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "10LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "`'0.02949"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "11BitrueCoinBTR"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "`'0.09274"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "12BitMartTokenBMX"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "`'0.001644"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "13BitForexTokenBF"
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D15").Value = "`'3.221"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "14MCDexMCB"
$ws.Range("B16").Value = "CoinExToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D16").Value = "`'0.04729"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "15CoinExTokenCET"
$ws.Range("B17").Value = "One"
$ws.Range("C17").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D17").Value = "`'0.0005895"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "16OneONE"
$ws.Range("D18").Value = "`'0.006496"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "`'0.005017"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "18HotbitTokenHTBBestin24h"
$ws.Range("D20").Value = "`'0.001048"
$ws.Range("D20").Style = "Normal"
$ws.Range("D23").Value = "`'3.855"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Value = "`'6.425"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Value = "`'2.122"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").Value = "`'0.1298"
$ws.Range("D27").Style = "Normal"
$ws.Range("D40").Value = "`'0.04082"
$ws.Range("D40").Style = "Normal"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "`'0.1041"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "40BKEXTokenBKK"
$ws.Range("D42").Value = "`'0.002972"
$ws.Range("D42").Style = "Normal"
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").Value = "`'0.003261"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "42KickTokenKICK"
$ws.Range("D44").Value = "`'0.008625"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "`'0.00005837"
$ws.Range("D45").Style = "Normal"
$ws.Range("D48").Value = "`'0.009993"
$ws.Range("D48").Style = "Normal"
